# Swap data between row 13 and row 14 for the columns that differ:
# A, B, D, E, F, G, H, Q, R
# (Columns C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
#  are identical between the two rows, so no visible change for those.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range("$col`13")
    $cell14 = $ws.Range("$col`14")

    $v13 = $cell13.Value()
    $v14 = $cell14.Value()

    $cell13.Value = $v14
    $cell14.Value = $v13
}
